$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.290.44"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").Value = "2.597.77"
$ws.Range("E3").Value = "  +7.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.48"
$ws.Range("E5").Value = "  +4.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.60"
$ws.Range("E6").Value = "  +4.99%  "
$ws.Range("E7").Value = "  +6.22%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +15.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.33"
$ws.Range("E10").Value = "  +12.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.47"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0843"
$ws.Range("E12").Value = "  +8.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.20"
$ws.Range("E13").Value = "  +16.35%  "
$ws.Range("D14").Value = "2.990.21"
$ws.Range("E14").Value = "  +7.18%  "
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "2.608.33"
$ws.Range("E16").Value = "  +7.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.918"
$ws.Range("E17").Value = "  +9.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.94"
$ws.Range("E18").Value = "  +6.05%  "
$ws.Range("D19").Value = "46.465.36"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("E20").Value = "  +7.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.94"
$ws.Range("E21").Value = "  +4.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.67"
$ws.Range("E22").Value = "  +8.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.90"
$ws.Range("E23").Value = "  +7.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "272.35"
$ws.Range("E24").Value = "  +12.88%  "
$ws.Range("E25").Value = "  +8.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "30.25"
$ws.Range("E26").Value = "  +42.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.17"
$ws.Range("E27").Value = "  +12.00%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.02"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.54"
$ws.Range("E30").Value = "  +8.85%  "
$ws.Range("E31").Value = "  +4.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "39.08"
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("E33").Value = "  +13.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.63"
$ws.Range("E34").Value = "  -3.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.85"
$ws.Range("E35").Value = "  +3.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0840"
$ws.Range("E36").Value = "  +9.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.19"
$ws.Range("E37").Value = "  +10.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "150.21"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("E39").Value = "  +8.32%  "
$ws.Range("E40").Value = "  +5.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.08"
$ws.Range("E41").Value = "  +44.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.15"
$ws.Range("E42").Value = "  +9.47%  "
$ws.Range("E43").Value = "  +10.58%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.60"
$ws.Range("E44").Value = "  +11.50%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.10"
$ws.Range("E45").Value = "  +8.47%  "
$ws.Range("D46").Value = "2.138.38"
$ws.Range("E46").Value = "  +6.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "93.81"
$ws.Range("E48").Value = "  +5.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.74"
$ws.Range("E49").Value = "  +13.83%  "
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "108.79"
$ws.Range("E51").Value = "  +7.87%  "
